$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new vocabulary rows (115-117) with a new "source" column D
$ws.Range("A115").Value = "Grande"
$ws.Range("B115").Value = "Toto"

$ws.Range("A116").Value = "Gran"
$ws.Range("B116").Value = "Toto"

$ws.Range("A117").Value = "Pequeño"
$ws.Range("B117").Value = "Mimi"

$ws.Range("D115").Value = "ChatGPT"
$ws.Range("D116").Value = "ChatGPT"
$ws.Range("D117").Value = "ChatGPT"

# Scroll the view down to show the new rows, matching the saved view state
$ws.Range("D118").Select()
$excel.ActiveWindow.ScrollRow = 97
